$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 24.1875
$ws.Range("I11").Value = 24.1875
$ws.Range("K11").Value = 24.1875
$ws.Range("M11").Value = 115.8125

$ws.Range("H17").Value = 3763.8333
$ws.Range("J17").Value = 3867.5881
$ws.Range("L17").Value = 11602.7643
$ws.Range("N17").Value = -11938.7643

$ws.Range("H40").Value = 2164.7693
$ws.Range("I40").Value = 2038.2778
$ws.Range("K40").Value = 2038.2778
$ws.Range("M40").Value = -1863.2778

$ws.Range("H51").Value = 10305.944
$ws.Range("I51").Value = 2999
$ws.Range("K51").Value = 2999
$ws.Range("M51").Value = -2515

$ws.Range("H98").Value = 2880.9333
$ws.Range("I98").Value = 775.75
$ws.Range("K98").Value = 775.75
$ws.Range("M98").Value = 722.25

$ws.Range("H122").Value = 2880.9333
$ws.Range("I122").Value = 775.75
$ws.Range("K122").Value = 2327.25
$ws.Range("M122").Value = 122.75

$ws.Range("H131").Value = 9218.615
$ws.Range("I131").Value = 1906
$ws.Range("K131").Value = 5718
$ws.Range("M131").Value = -678

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 930.875
$ws.Range("I2").Value = 706.7143
$ws.Range("K2").Value = 706.7143
$ws.Range("M2").Value = -593.7143

$ws.Range("H63").Value = 27699
$ws.Range("I63").Value = 41415
$ws.Range("K63").Value = 41415
$ws.Range("M63").Value = -40729

$ws.Range("H66").Value = 27699
$ws.Range("I66").Value = 41415
$ws.Range("K66").Value = 207075
$ws.Range("M66").Value = -203643

$ws.Range("H74").Value = 1767.0526
$ws.Range("I74").Value = 1874.909
$ws.Range("J74").Value = 1618.75
$ws.Range("K74").Value = 1874.909
$ws.Range("L74").Value = 1618.75
$ws.Range("M74").Value = -1000.909
$ws.Range("N74").Value = -3366.75

$ws.Range("H77").Value = 1767.0526
$ws.Range("I77").Value = 1874.909
$ws.Range("J77").Value = 1618.75
$ws.Range("K77").Value = 9374.545
$ws.Range("L77").Value = 8093.75
$ws.Range("M77").Value = -5006.545
$ws.Range("N77").Value = -16829.75

$ws.Range("H88").Value = 566.1667
$ws.Range("I88").Value = 559.4
$ws.Range("K88").Value = 559.4
$ws.Range("M88").Value = -153.4

$ws.Range("H91").Value = 566.1667
$ws.Range("I91").Value = 559.4
$ws.Range("K91").Value = 559.4
$ws.Range("M91").Value = 844.6

$ws.Range("H116").Value = 930.875
$ws.Range("I116").Value = 706.7143
$ws.Range("K116").Value = 706.7143
$ws.Range("M116").Value = 1587.2857

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 930.875
$ws.Range("I3").Value = 706.7143
$ws.Range("K3").Value = 706.7143
$ws.Range("M3").Value = -592.7143

$ws.Range("H99").Value = 2211.4546
$ws.Range("I99").Value = 1915.875
$ws.Range("K99").Value = 1915.875
$ws.Range("M99").Value = -417.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1997.3846
$ws.Range("I31").Value = 1522
$ws.Range("J31").Value = 2404.8572
$ws.Range("K31").Value = 1522
$ws.Range("L31").Value = 2404.8572
$ws.Range("M31").Value = -1227
$ws.Range("N31").Value = -2994.8572

$ws.Range("H34").Value = 1997.3846
$ws.Range("I34").Value = 1522
$ws.Range("J34").Value = 2404.8572
$ws.Range("K34").Value = 1522
$ws.Range("L34").Value = 2404.8572
$ws.Range("M34").Value = -1320
$ws.Range("N34").Value = -2808.8572

$ws.Range("H41").Value = 900
$ws.Range("I41").Value = 900
$ws.Range("K41").Value = 900
$ws.Range("M41").Value = -472

$ws.Range("H132").Value = 3149.7144
$ws.Range("I132").Value = 4250
$ws.Range("J132").Value = 2324.5
$ws.Range("K132").Value = 12750
$ws.Range("L132").Value = 6973.5
$ws.Range("M132").Value = -10220
$ws.Range("N132").Value = -12033.5

$ws.Range("H134").Value = 2343.6428
$ws.Range("I134").Value = 2609.2727
$ws.Range("J134").Value = 1369.6666
$ws.Range("K134").Value = 7827.8181
$ws.Range("L134").Value = 4108.9998
$ws.Range("M134").Value = -5292.8181
$ws.Range("N134").Value = -9178.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 79951
$ws.Range("J37").Value = 79951
$ws.Range("L37").Value = 239853
$ws.Range("N37").Value = -240077

$ws.Range("H51").Value = 1999.25
$ws.Range("J51").Value = 1999.25
$ws.Range("L51").Value = 5997.75
$ws.Range("N51").Value = -6917.75

$ws.Range("H80").Value = 1560
$ws.Range("I80").Value = 1560
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 4680
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -3744
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 1560
$ws.Range("I83").Value = 1560
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 14040
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -9360
$ws.Range("N83").ClearContents()

$ws.Range("H128").Value = 800000
$ws.Range("I128").Value = 800000
$ws.Range("K128").Value = 2400000
$ws.Range("M128").Value = -2395020

$ws.Range("H132").Value = 2229
$ws.Range("I132").Value = 959.3333
$ws.Range("K132").Value = 8633.9997
$ws.Range("M132").Value = -6103.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 461.92307
$ws.Range("J97").Value = 514
$ws.Range("L97").Value = 514
$ws.Range("N97").Value = -1506

$ws.Range("H113").Value = 3111.625
$ws.Range("J113").Value = 3798.8
$ws.Range("L113").Value = 3798.8
$ws.Range("N113").Value = -8138.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7637.0557
$ws.Range("I7").Value = 4492.3335
$ws.Range("K7").Value = 4492.3335
$ws.Range("M7").Value = -4380.3335

$ws.Range("H13").Value = 1000
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws.Range("H46").Value = 3000
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H61").Value = 2147.5
$ws.Range("I61").Value = 1860.2858
$ws.Range("K61").Value = 1860.2858
$ws.Range("M61").Value = -1658.2858

$ws.Range("H68").Value = 2428.1428
$ws.Range("I68").Value = 2199.4
$ws.Range("K68").Value = 2199.4
$ws.Range("M68").Value = -1450.4

$ws.Range("H71").Value = 2428.1428
$ws.Range("I71").Value = 2199.4
$ws.Range("K71").Value = 10997
$ws.Range("M71").Value = -7253

$ws.Range("H100").Value = 3710.4375
$ws.Range("I100").Value = 2988.0908
$ws.Range("K100").Value = 2988.0908
$ws.Range("M100").Value = -2447.0908

$ws.Range("H113").Value = 2147.5
$ws.Range("I113").Value = 1860.2858
$ws.Range("K113").Value = 1860.2858
$ws.Range("M113").Value = 309.7141999999999

$ws.Range("H126").Value = 7637.0557
$ws.Range("I126").Value = 4492.3335
$ws.Range("K126").Value = 13477.0005
$ws.Range("M126").Value = -11007.0005

$ws.Range("H132").Value = 2723.3076
$ws.Range("I132").Value = 2354.3635
$ws.Range("K132").Value = 7063.0905
$ws.Range("M132").Value = -4533.0905

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1994.6
$ws.Range("I107").Value = 1994.5
$ws.Range("K107").Value = 5983.5
$ws.Range("M107").Value = -4063.5

$ws.Range("H126").Value = 3749.182
$ws.Range("J126").Value = 6513.857
$ws.Range("L126").Value = 19541.571
$ws.Range("N126").Value = -24481.571

$ws.Range("H132").Value = 3554.652
$ws.Range("I132").Value = 3747.95
$ws.Range("J132").Value = 2266
$ws.Range("K132").Value = 11243.85
$ws.Range("L132").Value = 6798
$ws.Range("M132").Value = -8713.849999999999
$ws.Range("N132").Value = -11858
